$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.156477115925327
$ws.Range("C2").Value = 0.4295477836434998
$ws.Range("E2").Value = 0.06976766101828602
$ws.Range("F2").Value = 0.4443680307746121
$ws.Range("G2").Value = 0.002437069365432623
$ws.Range("I2").Value = 0.8319555317873011
$ws.Range("M2").Value = 0.4495498745744158
$ws.Range("N2").Value = 1.319102473191649
$ws.Range("B3").Value = 1.036782586654283
$ws.Range("C3").Value = 0.3795697224868491
$ws.Range("E3").Value = 0.06565594451598145
$ws.Range("F3").Value = 0.3878228170618172
$ws.Range("G3").Value = 0.002442718695768786
$ws.Range("I3").Value = 0.8181381695428556
$ws.Range("M3").Value = 0.4063453082006632
$ws.Range("N3").Value = 1.334082195588863
$ws.Range("B4").Value = 0.9638022806955746
$ws.Range("C4").Value = 0.3490538824271709
$ws.Range("E4").Value = 0.06318180856439071
$ws.Range("F4").Value = 0.3531389305169483
$ws.Range("G4").Value = 0.002446364628317675
$ws.Range("I4").Value = 0.8102741738912158
$ws.Range("M4").Value = 0.3800611554198738
$ws.Range("N4").Value = 1.34386441127991
$ws.Range("B5").Value = 0.9341884511616172
$ws.Range("C5").Value = 0.3366597044183948
$ws.Range("E5").Value = 0.06218606742718435
$ws.Range("F5").Value = 0.3390132514313251
$ws.Range("G5").Value = 0.002447895105149153
$ws.Range("I5").Value = 0.8072239722711529
$ws.Range("M5").Value = 0.3694100898623063
$ws.Range("N5").Value = 1.347996988388722
$ws.Range("B6").Value = 0.9292786448631887
$ws.Range("C6").Value = 0.3346041042978243
$ws.Range("E6").Value = 0.0620214738142586
$ws.Range("F6").Value = 0.336668177824194
$ws.Range("G6").Value = 0.002448151946080723
$ws.Range("I6").Value = 0.8067267744955728
$ws.Range("M6").Value = 0.3676450705304433
$ws.Range("N6").Value = 1.348692010428877
$ws.Range("B7").Value = 0.9634023905908862
$ws.Range("C7").Value = 0.3488865653177413
$ws.Range("E7").Value = 0.06316832934854943
$ws.Range("F7").Value = 0.3529483938368969
$ws.Range("G7").Value = 0.002446385087526599
$ws.Range("I7").Value = 0.8102324143754771
$ws.Range("M7").Value = 0.3799172704627125
$ws.Range("N7").Value = 1.343919553512222
$ws.Range("B8").Value = 1.115098394762697
$ws.Range("C8").Value = 0.4122788374519928
$ws.Range("E8").Value = 0.06833934100840011
$ws.Range("F8").Value = 0.4248636149813905
$ws.Range("G8").Value = 0.002438980572610739
$ws.Range("I8").Value = 0.8270617616444156
$ws.Range("M8").Value = 0.4346015617024577
$ws.Range("N8").Value = 1.324145689255879
$ws.Range("B9").Value = 1.416778345419914
$ws.Range("C9").Value = 0.5380301077658487
$ws.Range("E9").Value = 0.0788899659015172
$ws.Range("F9").Value = 0.5661985755042025
$ws.Range("G9").Value = 0.002425858901975469
$ws.Range("I9").Value = 0.8650510743996165
$ws.Range("M9").Value = 0.543836105923134
$ws.Range("N9").Value = 1.290041591060351
$ws.Range("B10").Value = 1.641200245777213
$ws.Range("C10").Value = 0.6314219508899441
$ws.Range("E10").Value = 0.08690612312955892
$ws.Range("F10").Value = 0.6702781546542269
$ws.Range("G10").Value = 0.00241706031992454
$ws.Range("I10").Value = 0.8961010096746236
$ws.Range("M10").Value = 0.6254092015546462
$ws.Range("N10").Value = 1.267879928230911
$ws.Range("B11").Value = 1.743947260240589
$ws.Range("C11").Value = 0.6741539053582528
$ws.Range("E11").Value = 0.09061356848916802
$ws.Range("F11").Value = 0.7176906081379002
$ws.Range("G11").Value = 0.002413238137418899
$ws.Range("I11").Value = 0.9109299847174555
$ws.Range("M11").Value = 0.6628276418996251
$ws.Range("N11").Value = 1.258436679861369
$ws.Range("B12").Value = 1.782952809628057
$ws.Range("C12").Value = 0.69037310133956
$ws.Range("E12").Value = 0.09202648310098738
$ws.Range("F12").Value = 0.7356546913071611
$ws.Range("G12").Value = 0.002411816532501682
$ws.Range("I12").Value = 0.9166483334424527
$ws.Range("M12").Value = 0.6770433570585652
$ws.Range("N12").Value = 1.254953419313459
$ws.Range("B13").Value = 1.774547880907335
$ws.Range("C13").Value = 0.686878306830863
$ws.Range("E13").Value = 0.09172178324661218
$ws.Range("F13").Value = 0.7317853510981394
$ws.Range("G13").Value = 0.002412121556698009
$ws.Range("I13").Value = 0.9154121796199632
$ws.Range("M13").Value = 0.6739796707122565
$ws.Range("N13").Value = 1.255699464455688
$ws.Range("B14").Value = 1.747154298519376
$ws.Range("C14").Value = 0.675487504953594
$ws.Range("E14").Value = 0.09072962838283161
$ws.Range("F14").Value = 0.7191683204515869
$ws.Range("G14").Value = 0.002413120665434212
$ws.Range("I14").Value = 0.9113983640067715
$ws.Range("M14").Value = 0.6639962433113027
$ws.Range("N14").Value = 1.258148246902131
$ws.Range("B15").Value = 1.730387726923198
$ws.Range("C15").Value = 0.668515255774139
$ws.Range("E15").Value = 0.09012308209408815
$ws.Range("F15").Value = 0.7114413442032514
$ws.Range("G15").Value = 0.002413736000032439
$ws.Range("I15").Value = 0.9089532416824397
$ws.Range("M15").Value = 0.6578871634081906
$ws.Range("N15").Value = 1.259660295720735
$ws.Range("B16").Value = 1.634498957983794
$ws.Range("C16").Value = 0.6286344630266285
$ws.Range("E16").Value = 0.08666507738066542
$ws.Range("F16").Value = 0.6671810134426437
$ws.Range("G16").Value = 0.002417313723517896
$ws.Range("I16").Value = 0.8951462141992863
$ws.Range("M16").Value = 0.6229701947742825
$ws.Range("N16").Value = 1.268509987298849
$ws.Range("B17").Value = 1.575844599763002
$ws.Range("C17").Value = 0.6042336558546708
$ws.Range("E17").Value = 0.0845594556157252
$ws.Range("F17").Value = 0.6400460337215605
$ws.Range("G17").Value = 0.002419554616135941
$ws.Range("I17").Value = 0.8868575968340338
$ws.Range("M17").Value = 0.6016303495644451
$ws.Range("N17").Value = 1.274103095486616
$ws.Range("B18").Value = 1.542169755060854
$ws.Range("C18").Value = 0.5902221383930737
$ws.Range("E18").Value = 0.08335407137427353
$ws.Range("F18").Value = 0.6244449056556647
$ws.Range("G18").Value = 0.002420860501774911
$ws.Range("I18").Value = 0.8821563847635758
$ws.Range("M18").Value = 0.5893853612398061
$ws.Range("N18").Value = 1.277380118236316
$ws.Range("B19").Value = 1.530778516183148
$ws.Range("C19").Value = 0.585482002315473
$ws.Range("E19").Value = 0.08294692323492114
$ws.Range("F19").Value = 0.6191636801734006
$ws.Range("G19").Value = 0.00242130557405931
$ws.Range("I19").Value = 0.8805759512397771
$ws.Range("M19").Value = 0.585244375076087
$ws.Range("N19").Value = 1.278499941632781
$ws.Range("B20").Value = 1.582082059670711
$ws.Range("C20").Value = 0.6068287503193801
$ws.Range("E20").Value = 0.0847830096140072
$ws.Range("F20").Value = 0.6429339538360921
$ws.Range("G20").Value = 0.002419314312710732
$ws.Range("I20").Value = 0.8877330727319048
$ws.Range("M20").Value = 0.603898987113169
$ws.Range("N20").Value = 1.273501481648488
$ws.Range("B21").Value = 1.755197786938993
$ws.Range("C21").Value = 0.6788322253649994
$ws.Range("E21").Value = 0.09102080256211309
$ws.Range("F21").Value = 0.7228739723492197
$ws.Range("G21").Value = 0.002412826504876744
$ws.Range("I21").Value = 0.9125745117824522
$ws.Range("M21").Value = 0.6669273527724755
$ws.Range("N21").Value = 1.257426457081934
$ws.Range("B22").Value = 1.868908626578843
$ws.Range("C22").Value = 0.7261103486205798
$ws.Range("E22").Value = 0.0951500471636777
$ws.Range("F22").Value = 0.7751780083420101
$ws.Range("G22").Value = 0.002408736495020099
$ws.Range("I22").Value = 0.9294106227636547
$ws.Range("M22").Value = 0.7083898402175066
$ws.Range("N22").Value = 1.24746122024883
$ws.Range("B23").Value = 1.808165841268249
$ws.Range("C23").Value = 0.7008563739688611
$ws.Range("E23").Value = 0.09294131050838672
$ws.Range("F23").Value = 0.7472568307916134
$ws.Range("G23").Value = 0.002410905724450796
$ws.Range("I23").Value = 0.9203693410757978
$ws.Range("M23").Value = 0.686235331659276
$ws.Range("N23").Value = 1.252730070970173
$ws.Range("B24").Value = 1.579261959980954
$ws.Range("C24").Value = 0.6056554558310836
$ws.Range("E24").Value = 0.0846819247789341
$ws.Range("F24").Value = 0.6416283278902171
$ws.Range("G24").Value = 0.002419422899122003
$ws.Range("I24").Value = 0.8873370707829196
$ws.Range("M24").Value = 0.6028732627677016
$ws.Range("N24").Value = 1.273773279867783
$ws.Range("B25").Value = 1.334693862315419
$ws.Range("C25").Value = 0.5038441710153734
$ws.Range("E25").Value = 0.07599033388886767
$ws.Range("F25").Value = 0.5279251897347308
$ws.Range("G25").Value = 0.002429260042435939
$ws.Range("I25").Value = 0.854229916618138
$ws.Range("M25").Value = 0.5140610463411832
$ws.Range("N25").Value = 1.298762715148428
